# Nexial #system sheet update:
#  - insert a new "macro" column (M), shifting mail..xml columns one to the right
#  - insert a new "macro" category row in the target (A) column
#  - correct/extend the "external" (H) column with two runProgram entries
#  - fix a typo and insert two new entries in the "web" (W, after shift) column
#  - update named ranges to reflect the new layout
#  - append the new function/category names to be used above

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1) Insert a brand new column at M ("macro"), shifting M:AA -> N:AB
# ---------------------------------------------------------------------------
$ws.Columns("M").Insert()

$ws.Range("M1").Value = "macro"
$ws.Range("M2").Value = "description()"
$ws.Range("M3").Value = "expects(var,default)"
$ws.Range("M4").Value = "produces(var,value)"

# ---------------------------------------------------------------------------
# 2) Insert a new row entry in the "target" (A) column for the new "macro"
#    category, between "json" (row 12) and the old "mail" (row 13)
# ---------------------------------------------------------------------------
$ws.Range("A13").Insert(-4121)   # xlShiftDown
$ws.Range("A13").Value = "macro"

# ---------------------------------------------------------------------------
# 3) "external" (H) column: fix the spelling of the existing runProgram entry
#    and append the new runProgramNoWait entry (no shifting needed - appended
#    at the end of the existing list)
# ---------------------------------------------------------------------------
$ws.Range("H3").Value = "runProgram(programPathAndParams)"
$ws.Range("H4").Value = "runProgramNoWait(programPathAndParams)"

# ---------------------------------------------------------------------------
# 4) "web" (now column W) column updates:
#    a) fix typo assertIENavtiveMode() -> assertIENativeMode() in place
#    b) insert clickOffset(locator,x,y) after clickByLabelAndWait
#    c) insert saveAttributeList(var,locator,attrName) after saveAttribute
# ---------------------------------------------------------------------------
$ws.Range("W21").Value = "assertIENativeMode()"

$ws.Range("W51").Insert(-4121)   # xlShiftDown
$ws.Range("W51").Value = "clickOffset(locator,x,y)"

$ws.Range("W82").Insert(-4121)   # xlShiftDown
$ws.Range("W82").Value = "saveAttributeList(var,locator,attrName)"

# ---------------------------------------------------------------------------
# 5) Update the named ranges so they refer to the correct (shifted) columns
#    and row extents
# ---------------------------------------------------------------------------
$wb.Names.Add("macro", "='#system'!`$M`$2:`$M`$4")

$wb.Names.Item("external").RefersTo  = "='#system'!`$H`$2:`$H`$4"
$wb.Names.Item("mail").RefersTo      = "='#system'!`$N`$2:`$N`$2"
$wb.Names.Item("number").RefersTo    = "='#system'!`$O`$2:`$O`$15"
$wb.Names.Item("pdf").RefersTo       = "='#system'!`$P`$2:`$P`$16"
$wb.Names.Item("rdbms").RefersTo     = "='#system'!`$Q`$2:`$Q`$7"
$wb.Names.Item("redis").RefersTo     = "='#system'!`$R`$2:`$R`$10"
$wb.Names.Item("sms").RefersTo       = "='#system'!`$S`$2:`$S`$2"
$wb.Names.Item("sound").RefersTo     = "='#system'!`$T`$2:`$T`$5"
$wb.Names.Item("ssh").RefersTo       = "='#system'!`$U`$2:`$U`$9"
$wb.Names.Item("step").RefersTo      = "='#system'!`$V`$2:`$V`$4"
$wb.Names.Item("target").RefersTo    = "='#system'!`$A`$2:`$A`$28"
$wb.Names.Item("web").RefersTo       = "='#system'!`$W`$2:`$W`$122"
$wb.Names.Item("webalert").RefersTo  = "='#system'!`$X`$2:`$X`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$Y`$2:`$Y`$8"
$wb.Names.Item("ws").RefersTo        = "='#system'!`$Z`$2:`$Z`$17"
$wb.Names.Item("ws.async").RefersTo  = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("xml").RefersTo       = "='#system'!`$AB`$2:`$AB`$13"
